$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells: prepend "Geography " to the existing headers.
$ws.Range("A1").Value = "Geography Code"
$ws.Range("B1").Value = "Geography Name"
$ws.Range("C1").Value = "Geography Colour Code"

# Update the structured table's column names to match, if a ListObject is present.
foreach ($lo in $ws.ListObjects) {
    $lo.ListColumns.Item(1).Name = "Geography Code"
    $lo.ListColumns.Item(2).Name = "Geography Name"
    $lo.ListColumns.Item(3).Name = "Geography Colour Code"
}

# Update the selected cell to C2 (matching the saved selection in the file).
$ws.Range("C2").Select()
